# Weekly update: insert a new weekly price record for Jengibre (Vega Central
# Mapocho de Santiago) as a new row 31, pushing the existing rows 31-63 down
# to 32-64 (dimension grows from A1:R63 to A1:R64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31 - shifts old rows 31..63 down to 32..64.
$ws.Rows(31).Insert()

# Fill in the new row's data.
$ws.Cells(31, 1).Value  = 9
$ws.Cells(31, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells(31, 3).Value  = "Metropolitana"
$ws.Cells(31, 4).Value  = 44494
$ws.Cells(31, 5).Value  = 13
$ws.Cells(31, 6).Value  = 100114007
$ws.Cells(31, 7).Value  = "Jengibre"
$ws.Cells(31, 8).Value  = "Sin especificar"
$ws.Cells(31, 9).Value  = "Primera"
$ws.Cells(31, 10).Value = 780
$ws.Cells(31, 11).Value = 15000
$ws.Cells(31, 12).Value = 15000
$ws.Cells(31, 13).Value = 15000
$ws.Cells(31, 14).Value = "$/caja 13 kilos"
$ws.Cells(31, 15).Value = "Perú"
$ws.Cells(31, 16).Value = 1154
$ws.Cells(31, 17).Value = 13
$ws.Cells(31, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other rows' "Fecha" column.
$ws.Cells(31, 4).NumberFormat = $ws.Cells(32, 4).NumberFormat
